$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.295.18"
$ws.Range("E2").Value = "  +6.87%  "
$ws.Range("D3").Value = "3.683.53"
$ws.Range("E3").Value = "  +19.19%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'598.28"
$ws.Range("E5").Value = "  +3.73%  "
$ws.Range("D6").Value = "'184.08"
$ws.Range("E6").Value = "  +6.51%  "
$ws.Range("D7").Value = "3.681.07"
$ws.Range("E7").Value = "  +19.16%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  +4.52%  "
$ws.Range("E10").Value = "  +7.54%  "
$ws.Range("D11").Value = "'6.56"
$ws.Range("E11").Value = "  +3.24%  "
$ws.Range("D12").Value = "'0.500"
$ws.Range("E12").Value = "  +6.52%  "
$ws.Range("D13").Value = "'40.00"
$ws.Range("E13").Value = "  +11.00%  "
$ws.Range("D14").Value = "'0.0000253"
$ws.Range("E14").Value = "  +6.04%  "
$ws.Range("D15").Value = "4.291.82"
$ws.Range("E15").Value = "  +19.14%  "
$ws.Range("D16").Value = "71.240.11"
$ws.Range("E16").Value = "  +6.91%  "
$ws.Range("D17").Value = "3.671.28"
$ws.Range("E17").Value = "  +18.88%  "
$ws.Range("E18").Value = "  +1.99%  "
$ws.Range("D19").Value = "'7.50"
$ws.Range("E19").Value = "  +7.80%  "
$ws.Range("D20").Value = "'16.97"
$ws.Range("E20").Value = "  +1.05%  "
$ws.Range("D21").Value = "'518.01"
$ws.Range("E21").Value = "  +6.43%  "
$ws.Range("D22").Value = "'9.21"
$ws.Range("E22").Value = "  +17.54%  "
$ws.Range("D23").Value = "'0.744"
$ws.Range("E23").Value = "  +8.42%  "
$ws.Range("D24").Value = "'87.79"
$ws.Range("E24").Value = "  +5.27%  "
$ws.Range("D25").Value = "'13.54"
$ws.Range("E25").Value = "  +6.34%  "
$ws.Range("E26").Value = "  +8.08%  "
$ws.Range("D27").Value = "'10.80"
$ws.Range("E27").Value = "  +7.21%  "
$ws.Range("E28").Value = "  -0.02%  "
$ws.Range("D29").Value = "'2.52"
$ws.Range("E29").Value = "  +12.28%  "
$ws.Range("E30").Value = "  +2.30%  "
$ws.Range("D31").Value = "'31.91"
$ws.Range("E31").Value = "  +13.92%  "
$ws.Range("E32").Value = "  +6.80%  "
$ws.Range("E33").Value = "  +17.24%  "
$ws.Range("E34").Value = "  +4.13%  "
$ws.Range("E35").Value = "  -0.07%  "
$ws.Range("D36").Value = "'6.19"
$ws.Range("E36").Value = "  +10.46%  "
$ws.Range("E37").Value = "  +8.01%  "
$ws.Range("E38").Value = "  +11.72%  "
$ws.Range("E39").Value = "  +8.44%  "
$ws.Range("D40").Value = "'50.75"
$ws.Range("E40").Value = "  +3.64%  "
$ws.Range("D41").Value = "'46.31"
$ws.Range("E42").Value = "  +3.95%  "
$ws.Range("D43").Value = "3.188.40"
$ws.Range("E43").Value = "  +14.74%  "
$ws.Range("D44").Value = "'8.80"
$ws.Range("E44").Value = "  +6.88%  "
$ws.Range("D45").Value = "'2.75"
$ws.Range("E45").Value = "  +6.07%  "
$ws.Range("D46").Value = "'400.36"
$ws.Range("E46").Value = "  +8.96%  "
$ws.Range("E47").Value = "  +6.58%  "
$ws.Range("D48").Value = "'28.20"
$ws.Range("E48").Value = "  +15.25%  "
$ws.Range("D49").Value = "'135.84"
$ws.Range("E49").Value = "  +1.15%  "
$ws.Range("E50").Value = "  +0.00%  "
$ws.Range("D51").Value = "'2.43"
$ws.Range("E51").Value = "  +11.87%  "
